# Update test fixture for the "chunk reader" test: add a second real date
# (2023-02-20, serial 44977) below the existing date (2022-12-01, serial
# 44896), formatted as an ISO "yyyy-mm-dd" date, then leave the selection
# on the next empty cell (A3) as Excel would after typing/entering data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 44977
$ws.Range("A2").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("A3").Select() | Out-Null
